$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B7"   = 5.504999999999999
    "A9"   = -21.7984
    "B12"  = 5.402099999999995
    "B14"  = 6.1629
    "A18"  = -22.18850000000001
    "A20"  = -19.29319999999999
    "B26"  = 3.832100000000005
    "A27"  = -21.72569999999997
    "B27"  = 5.017600000000004
    "B29"  = 5.069399999999996
    "A35"  = -18.94509999999999
    "B37"  = 8.592100000000006
    "B38"  = 4.2193
    "B51"  = 5.804500000000002
    "B52"  = 5.379999999999999
    "B55"  = 5.211899999999998
    "A69"  = -21.6999
    "B69"  = 5.425099999999993
    "B70"  = 7.344400000000006
    "A76"  = -19.55009999999999
    "A78"  = -19.79279999999999
    "B81"  = 5.767
    "A82"  = -21.9186
    "A83"  = -21.71069999999999
    "B83"  = 5.543000000000004
    "A93"  = -20.80259999999997
    "B102" = 8.434000000000006
}

foreach ($cellRef in $changes.Keys) {
    $ws.Range($cellRef).Value = $changes[$cellRef]
}
